$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# The "CasesTab" query (row 2, column B) previously returned an extra
# `Cohort` column sourced from the (optional) cohort node. That column is
# no longer wanted, so drop the trailing RETURN clause line for it and
# remove the now-dangling trailing comma on the previous line.
$oldText = $ws.Range("B2").Value2
$cohortSuffix = ",`n        coalesce(co.cohort_description, '') AS ``Cohort``"
$newText = $oldText.Replace($cohortSuffix, "")
$ws.Range("B2").Value2 = $newText

# Reflect the edit in the active selection / view: the user last touched
# B2 (instead of B4) and the sheet is scrolled back to show row 1.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B2").Select()

$wb.Save()
